$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workload")

# Student identification row (was placeholder "<student id>" text)
$ws.Range("D5").Value = 5752132
$ws.Range("E5").Value = 5746957
$ws.Range("F5").Value = 5716934

# Student name row (was placeholder "<student name>" text)
$ws.Range("D6").Value = "Christiaan Baraya"
$ws.Range("E6").Value = "Ewout de Dobbelaar"
$ws.Range("F6").Value = "Egemen Yildiz"

# Basic features completion percentages
$ws.Range("E8").Value = 100
$ws.Range("D9").Value = 100
$ws.Range("D10").Value = 100
$ws.Range("F11").Value = 100
$ws.Range("D12").Value = 100
$ws.Range("D13").Value = 100
$ws.Range("F14").Value = 100
$ws.Range("D15").Value = 100

# Extra features completion percentages
$ws.Range("D19").Value = 100
$ws.Range("D20").Value = 100
$ws.Range("F21").Value = 100
$ws.Range("E22").Value = 100
$ws.Range("D23").Value = 100
$ws.Range("E24").Value = 100

# View: zoom level and active selection as left by the author
$excel.ActiveWindow.Zoom = 71
$ws.Range("H15").Select()

$wb.Save()
